$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.512.25'
$ws.Range("E2").Value = '  -4.11%  '
$ws.Range("D3").Value = '2.959.30'
$ws.Range("E3").Value = '  -6.05%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.22'
$ws.Range("E5").Value = '  -4.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.34'
$ws.Range("E6").Value = '  -6.88%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.570'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '2.966.41'
$ws.Range("E9").Value = '  -5.90%  '
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.12'
$ws.Range("E11").Value = '  -7.05%  '
$ws.Range("E12").Value = '  -3.33%  '
$ws.Range("D13").Value = '3.476.31'
$ws.Range("E13").Value = '  -6.27%  '
$ws.Range("E14").Value = '  -2.89%  '
$ws.Range("D15").Value = '61.577.71'
$ws.Range("E15").Value = '  -4.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.65'
$ws.Range("E16").Value = '  -5.38%  '
$ws.Range("D17").Value = '2.966.67'
$ws.Range("E17").Value = '  -5.70%  '
$ws.Range("E18").Value = '  -4.83%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '382.71'
$ws.Range("E20").Value = '  -4.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.99'
$ws.Range("E21").Value = '  -4.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.68'
$ws.Range("E22").Value = '  -5.61%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.39'
$ws.Range("E24").Value = '  -4.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.471'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = '3.083.60'
$ws.Range("E27").Value = '  -3.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").Value = '0.0₃0932'
$ws.Range("E29").Value = '  -7.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.34'
$ws.Range("E30").Value = '  -5.01%  '
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.45'
$ws.Range("E33").Value = '  -3.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '159.47'
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.95'
$ws.Range("E36").Value = '  -4.89%  '
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("E38").Value = '  -4.35%  '
$ws.Range("E40").Value = '  -3.15%  '
$ws.Range("D41").Value = '2.407.93'
$ws.Range("E41").Value = '  -9.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.21'
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.10'
$ws.Range("E43").Value = '  -6.42%  '
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.95'
$ws.Range("E48").Value = '  -8.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0958'
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.76'
$ws.Range("E50").Value = '  -6.24%  '
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '267.95'
$ws.Range("E51").Value = '  -6.79%  '
